# Refresh the cryptos price table (Sheet1) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are stored as plain text in this sheet
# (values like "1.01" or "543.18" would otherwise be auto-detected as numbers),
# so force a Text number format on every cell we are about to rewrite.
$priceRows = @(2,3,4,5,6,7,8,9,11,12,13,14,15,16,18,19,24,25,26,27,29,33,35,36,39,40,41,42,44,46,47,49)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$volumeRows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20,21,22,24,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $volumeRows) {
    $ws.Range("E$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '57.577.06'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '2.334.43'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.88%  '
$ws.Range("D5").Value = '543.18'
$ws.Range("E5").Value = '  +5.50%  '
$ws.Range("D6").Value = '135.56'
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '0.537'
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '2.371.27'
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("D12").Value = '0.153'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +4.33%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '23.75'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.760.75'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = '58.137.66'
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '2.357.19'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").Value = '339.88'
$ws.Range("E19").Value = '  +4.30%  '
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("E22").Value = '  +3.57%  '
$ws.Range("D24").Value = '62.42'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").Value = '0.170'
$ws.Range("D26").Value = '8.60'
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  +6.53%  '
$ws.Range("D29").Value = '174.94'
$ws.Range("E29").Value = '  +3.87%  '
$ws.Range("E30").Value = '  +4.81%  '
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = '18.60'
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("E34").Value = '  +13.19%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("E38").Value = '  +3.12%  '
$ws.Range("D39").Value = '1.61'
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("D40").Value = '39.44'
$ws.Range("E40").Value = '  +2.30%  '
$ws.Range("D41").Value = '150.53'
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '0.379'
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("D44").Value = '283.85'
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").Value = '0.0505'
$ws.Range("E46").Value = '  +1.18%  '
$ws.Range("D47").Value = '18.99'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("D49").Value = '0.0218'
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("E51").Value = '  +9.33%  '
